# Regenerate merged AHB files: rename the "_old"/"_new" diff-column header
# suffixes to the explicit format-version tags, freeze the header row, and
# wrap the used range in a native Excel table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2310", "_new" -> "_FV2404" ----
$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i] + "_FV2310"
}

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headers[$i] + "_FV2404"
}

# --- 2. Freeze the header row (split below row 1, top-left cell A2) ------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a native table ---------------------------
$rng = $ws.Range("A1:U71")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""
